$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.30544054395605258,
    -0.0099999999522282224,
    -0.0089999999514969176,
    -0.011999999986109344,
    -0.0059999999526896275,
    -0.0059999999515802926,
    -0.019999999941093805,
    -0.019999999940728763,
    0.037651242462737322,
    -0.0059999999499922296,
    -0.0044999999509656163,
    -0.0059999999496214151,
    -0.0059999999487985178,
    -0.011999999944128703,
    -0.0059999999483482114,
    -0.0059999999481614719,
    -0.005999999947913004,
    -0.0089999999456455981,
    -0.0089999999529206676,
    -0.0089999999524756902,
    -0.008999999952385096,
    -0.0089999999523433516,
    -0.0089999999509675632,
    -0.059889369195603948,
    -0.041999999925202758,
    -0.0059999999514168678,
    -0.0059999999513418167,
    -0.034339889595665696,
    -0.011999999946544548,
    -0.019999999940716329,
    -0.014999999944825149,
    -0.016199229919258151,
    0.066839425766595184
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
